$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 180, shifting existing rows (180-207) down to (181-208)
$ws.Rows(180).Insert()

# Fill in the new row 180 with data (copy of surrounding static columns + new values)
$ws.Cells.Item(180, 1).Value = 2
$ws.Cells.Item(180, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(180, 3).Value = "Coquimbo"
$ws.Cells.Item(180, 4).Value = 45015
$ws.Cells.Item(180, 5).Value = 4
$ws.Cells.Item(180, 6).Value = 100112043
$ws.Cells.Item(180, 7).Value = "Pepino ensalada"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 600
$ws.Cells.Item(180, 11).Value = 6500
$ws.Cells.Item(180, 12).Value = 7000
$ws.Cells.Item(180, 13).Value = 6750
$ws.Cells.Item(180, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(180, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(180, 16).Value = 96
$ws.Cells.Item(180, 17).Value = 70
$ws.Cells.Item(180, 18).Value = "Hortaliza"
